$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.556.47'
$ws.Range('E2').Value = '  +1.00%  '
$ws.Range('D3').Value = '2.644.01'
$ws.Range('E3').Value = '  +2.27%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '590.94'
$ws.Range('E5').Value = '  +1.46%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '145.54'
$ws.Range('E6').Value = '  -1.02%  '
$ws.Range('E8').Value = '  -1.00%  '
$ws.Range('D9').Value = '2.643.34'
$ws.Range('E9').Value = '  +2.28%  '
$ws.Range('E10').Value = '  -0.63%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '5.66'
$ws.Range('E11').Value = '  +0.00%  '
$ws.Range('E12').Value = '  +0.28%  '
$ws.Range('E13').Value = '  +0.17%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '27.58'
$ws.Range('E14').Value = '  +0.69%  '
$ws.Range('D15').Value = '3.114.13'
$ws.Range('E15').Value = '  +2.18%  '
$ws.Range('D16').Value = '63.451.68'
$ws.Range('E16').Value = '  +1.04%  '
$ws.Range('E17').Value = '  -0.64%  '
$ws.Range('D18').Value = '2.667.65'
$ws.Range('E18').Value = '  +2.97%  '
$ws.Range('E19').Value = '  +0.08%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '341.69'
$ws.Range('E20').Value = '  -0.24%  '
$ws.Range('E21').Value = '  -0.30%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.71'
$ws.Range('E22').Value = '  +0.27%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.69'
$ws.Range('E25').Value = '  +5.96%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.56'
$ws.Range('E26').Value = '  +7.75%  '
$ws.Range('E27').Value = '  -0.36%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '551.71'
$ws.Range('E28').Value = '  +17.60%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.48'
$ws.Range('E29').Value = '  +1.72%  '
$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  +0.11%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.80'
$ws.Range('E31').Value = '  -0.26%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.81'
$ws.Range('E32').Value = '  +12.93%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.97'
$ws.Range('E33').Value = '  +2.47%  '
$ws.Range('D34').Value = '0.0₃0808'
$ws.Range('E34').Value = '  -1.65%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '175.15'
$ws.Range('E35').Value = '  -0.53%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.88'
$ws.Range('E36').Value = '  +8.14%  '
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.399'
$ws.Range('E38').Value = '  -0.86%  '
$ws.Range('E39').Value = '  +0.43%  '
$ws.Range('E40').Value = '  +4.12%  '
$ws.Range('E41').Value = '  +0.09%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '170.03'
$ws.Range('E42').Value = '  +7.91%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '40.36'
$ws.Range('E43').Value = '  +2.41%  '
$ws.Range('E44').Value = '  -0.56%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '22.04'
$ws.Range('E45').Value = '  +4.07%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.630'
$ws.Range('E46').Value = '  -0.70%  '
$ws.Range('E47').Value = '  +2.09%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0959'
$ws.Range('E48').Value = '  -0.73%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0239'
$ws.Range('E49').Value = '  +1.45%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '18.77'
$ws.Range('E50').Value = '  +2.28%  '
$ws.Range('E51').Value = '  -1.06%  '
